# ------------------------------------------------------------------
# Update B083GT6ZBX_sales_po_comparison workbook with the corrected
# forecast output: rename/re-layout the main sheet (insert an
# "Order Week" column) and add three new report sheets
# (Weekly Growth, Volume Insights, Prediction Info).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet --------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# --- 2. Insert the new "Order Week" column (C) ---------------------
# Shift the old PO_Requested_Qty column from C to D, and put a new
# "Order Week" date column in C.
$ws1.Columns.Item(3).Insert()

# --- 3. Header row ---------------------------------------------------
$ws1.Range("A1").Value = "ds"
$ws1.Range("B1").Value = "y"
$ws1.Range("C1").Value = "Order Week"
$ws1.Range("D1").Value = "PO_Requested_Qty"

# Re-apply the bold/border/centered header style to the newly
# inserted C1 and to D1 (reuses A1's existing style).
$ws1.Range("A1").Copy()
$ws1.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Body data (rows 2..103) --------------------------------------
# Columns: A=ds (week-ending sales date), B=y (sales qty),
#          C=Order Week (PO date), D=PO_Requested_Qty
$raw = "44927,9,44921,0;44927,0,44921,0;44941,1,44935,0;44941,1,44935,0;44955,0,44949,0;44962,0,44956,0;44976,14,44970,0;44990,0,44984,0;44997,0,44991,0;45004,0,44998,0;45011,0,45005,0;45018,0,45012,0;45025,1,45019,0;45032,0,45026,0;45039,0,45033,0;45046,0,45040,0;45053,0,45047,0;45060,1,45054,0;45067,0,45061,0;45074,0,45068,0;45081,0,45075,0;45088,0,45082,0;45095,0,45089,0;45102,0,45096,0;45109,0,45103,0;45116,0,45110,0;45123,0,45117,0;45130,0,45124,0;45137,0,45131,0;45144,0,45138,0;45151,0,45145,0;45158,0,45152,0;45165,0,45159,0;45172,0,45166,0;45179,0,45173,0;45186,0,45180,0;45193,0,45187,0;45200,0,45194,0;45207,0,45201,0;45214,0,45208,0;45221,0,45215,0;45228,0,45222,0;45235,0,45229,0;45242,0,45236,0;45249,0,45243,0;45256,0,45250,0;45263,0,45257,0;45270,0,45264,0;45277,0,45271,0;45291,0,45285,0;45298,0,45292,0;45305,0,45299,0;45312,0,45306,0;45319,0,45313,0;45326,0,45320,0;45333,0,45327,0;45340,0,45334,0;45347,0,45341,0;45354,0,45348,0;45361,0,45355,0;45368,0,45362,0;45375,0,45369,0;45382,0,45376,0;45396,1,45390,0;45403,0,45397,0;45410,0,45404,0;45417,0,45411,0;45417,0,45411,0;45424,0,45418,0;45431,0,45425,0;45438,0,45432,0;45445,0,45439,0;45452,0,45446,0;45459,0,45453,0;45466,0,45460,0;45473,0,45467,0;45480,0,45474,0;45487,0,45481,0;45494,0,45488,0;45501,0,45495,0;45508,0,45502,0;45515,0,45509,0;45522,0,45516,0;45529,1,45523,0;45536,0,45530,0;45543,0,45537,0;45550,0,45544,0;45557,0,45551,0;45564,0,45558,0;45571,0,45565,0;45578,0,45572,0;45585,0,45579,0;45592,0,45586,0;45599,0,45593,0;45606,0,45600,0;45613,0,45607,0;45620,0,45614,0;45627,0,45621,0;45634,0,45628,0;45641,0,45635,0;45648,0,45642,0;45655,0,45649,0"

$recordStrings = $raw.Split(";")
$rowCount = $recordStrings.Length

$values = New-Object 'object[,]' $rowCount,4
for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $recordStrings[$i].Split(",")
    $values[$i,0] = [double]$parts[0]
    $values[$i,1] = [double]$parts[1]
    $values[$i,2] = [double]$parts[2]
    $values[$i,3] = [double]$parts[3]
}

$lastRow = 1 + $rowCount
$ws1.Range("A2:D$lastRow").Value = $values

# --- 5. Re-apply the date number format / style to columns A and C ---
$ws1.Range("A2").Copy()
$ws1.Range("A2:A$lastRow").PasteSpecial(-4122)
$ws1.Range("C2:C$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Add the "Weekly Growth" sheet --------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"
$ws1.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 7. Add the "Volume Insights" sheet ------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"
$ws1.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0

# --- 8. Add the "Prediction Info" sheet ------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("A2").Value = 0

# --- 9. Restore the active sheet/selection ---------------------------
$ws1.Activate()
$ws1.Range("A1").Select()

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
